$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2550.4
$ws.Range("J17").Value = 2550.4
$ws.Range("L17").Value = 7651.200000000001
$ws.Range("N17").Value = -7987.200000000001

# Row 40
$ws.Range("H40").Value = 2841.5715

# Row 41
$ws.Range("H41").Value = 2436.9285
$ws.Range("I41").Value = 2744
$ws.Range("K41").Value = 2744
$ws.Range("M41").Value = -2304

# Row 62
$ws.Range("H62").Value = 1746.75
$ws.Range("I62").Value = 1746.75
$ws.Range("K62").Value = 1746.75
$ws.Range("M62").Value = -1122.75

# Row 65
$ws.Range("H65").Value = 1746.75
$ws.Range("I65").Value = 1746.75
$ws.Range("K65").Value = 8733.75
$ws.Range("M65").Value = -5613.75

# Row 100
$ws.Range("H100").Value = 10118.3
$ws.Range("I100").Value = 1657.75
$ws.Range("J100").Value = 15758.667
$ws.Range("K100").Value = 1657.75
$ws.Range("L100").Value = 15758.667
$ws.Range("M100").Value = -1116.75
$ws.Range("N100").Value = -16840.667

# Row 116
$ws.Range("H116").Value = 3944.5
$ws.Range("I116").Value = 3944.5
$ws.Range("K116").Value = 3944.5
$ws.Range("M116").Value = -502.5

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 3073.3572
$ws.Range("I74").Value = 2488.3333
$ws.Range("J74").Value = 4126.4
$ws.Range("K74").Value = 2488.3333
$ws.Range("L74").Value = 4126.4
$ws.Range("M74").Value = -1614.3333
$ws.Range("N74").Value = -5874.4

# Row 77
$ws.Range("H77").Value = 3073.3572
$ws.Range("I77").Value = 2488.3333
$ws.Range("J77").Value = 4126.4
$ws.Range("K77").Value = 12441.6665
$ws.Range("L77").Value = 20632
$ws.Range("M77").Value = -8073.666499999999
$ws.Range("N77").Value = -29368

# Row 97
$ws.Range("H97").Value = 2276.5833
$ws.Range("I97").Value = 1496.8572
$ws.Range("K97").Value = 1496.8572
$ws.Range("M97").Value = -1000.8572

$ws = $wb.Worksheets.Item("BSM")
# Row 70
$ws.Range("H70").Value = 176464.67
$ws.Range("J70").Value = 176464.67
$ws.Range("L70").Value = 176464.67
$ws.Range("N70").Value = -177050.67

# Row 73
$ws.Range("H73").Value = 176464.67
$ws.Range("J73").Value = 176464.67
$ws.Range("L73").Value = 176464.67
$ws.Range("N73").Value = -178492.67

# Row 94
$ws.Range("H94").Value = 37970.285
$ws.Range("I94").Value = 13158.4
$ws.Range("K94").Value = 13158.4
$ws.Range("M94").Value = -12707.4

# Row 99
$ws.Range("H99").Value = 2705.5
$ws.Range("I99").Value = 1336
$ws.Range("K99").Value = 1336
$ws.Range("M99").Value = 162

# Row 105
$ws.Range("H105").Value = 3725.7058
$ws.Range("I105").Value = 3239.9092
$ws.Range("J105").Value = 4616.3335
$ws.Range("K105").Value = 3239.9092
$ws.Range("L105").Value = 4616.3335
$ws.Range("M105").Value = -1492.9092
$ws.Range("N105").Value = -8110.3335

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 617.6
$ws.Range("I22").Value = 322
$ws.Range("K22").Value = 322
$ws.Range("M22").Value = 28

# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 62
$ws.Range("H62").Value = 19192
$ws.Range("I62").Value = 7399
$ws.Range("J62").Value = 38847
$ws.Range("K62").Value = 7399
$ws.Range("L62").Value = 38847
$ws.Range("M62").Value = -6775
$ws.Range("N62").Value = -40095

# Row 65
$ws.Range("H65").Value = 19192
$ws.Range("I65").Value = 7399
$ws.Range("J65").Value = 38847
$ws.Range("K65").Value = 36995
$ws.Range("L65").Value = 194235
$ws.Range("M65").Value = -33875
$ws.Range("N65").Value = -200475

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 460.8889
$ws.Range("J23").Value = 514.6667
$ws.Range("L23").Value = 1544.0001
$ws.Range("N23").Value = -2014.0001

# Row 87
$ws.Range("H87").Value = 14862.429
$ws.Range("I87").Value = 14679
$ws.Range("K87").Value = 44037
$ws.Range("M87").Value = -42789

# Row 90
$ws.Range("H90").Value = 14862.429
$ws.Range("I90").Value = 14679
$ws.Range("K90").Value = 132111
$ws.Range("M90").Value = -125871

# Row 131
$ws.Range("H131").Value = 1441.32
$ws.Range("I131").Value = 965.38464
$ws.Range("J131").Value = 1956.9166
$ws.Range("K131").Value = 2896.15392
$ws.Range("L131").Value = 5870.7498
$ws.Range("M131").Value = 2143.84608
$ws.Range("N131").Value = -15950.7498

# Row 132
$ws.Range("H132").Value = 1101
$ws.Range("J132").Value = 2590
$ws.Range("L132").Value = 23310
$ws.Range("N132").Value = -28370

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1634.7778
$ws.Range("I132").Value = 1117.4615
$ws.Range("J132").Value = 2979.8
$ws.Range("K132").Value = 3352.3845
$ws.Range("L132").Value = 8939.400000000001
$ws.Range("M132").Value = -822.3844999999997
$ws.Range("N132").Value = -13999.4

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1387.6552
$ws.Range("J22").Value = 1416.5416
$ws.Range("L22").Value = 1416.5416
$ws.Range("N22").Value = -2006.5416

# Row 27
$ws.Range("H27").Value = 1387.6552
$ws.Range("J27").Value = 1416.5416
$ws.Range("L27").Value = 1416.5416
$ws.Range("N27").Value = -1630.5416

# Row 46
$ws.Range("H46").Value = 1707.5264
$ws.Range("J46").Value = 1843.1333
$ws.Range("L46").Value = 1843.1333
$ws.Range("N46").Value = -2219.1333

# Row 132
$ws.Range("H132").Value = 3132.45
$ws.Range("I132").Value = 2679.7144
$ws.Range("K132").Value = 8039.1432
$ws.Range("M132").Value = -5509.1432

# Row 136
$ws.Range("H136").Value = 2135.4688
$ws.Range("I136").Value = 1961
$ws.Range("J136").Value = 2426.25
$ws.Range("K136").Value = 5883
$ws.Range("L136").Value = 7278.75
$ws.Range("M136").Value = -3333
$ws.Range("N136").Value = -12378.75

$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 44535
$ws.Range("I51").Value = 39070
$ws.Range("K51").Value = 39070
$ws.Range("M51").Value = -38560

# Row 61
$ws.Range("H61").Value = 13360.333
$ws.Range("I61").Value = 10040.5
$ws.Range("K61").Value = 10040.5
$ws.Range("M61").Value = -9748.5

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 132
$ws.Range("H132").Value = 3463.5833
$ws.Range("I132").Value = 3922.75
$ws.Range("K132").Value = 11768.25
$ws.Range("M132").Value = -9238.25
